$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 799.0625
$ws.Range("J17").Value = 771.65515
$ws.Range("L17").Value = 2314.96545
$ws.Range("N17").Value = -2650.96545

$ws.Range("H18").Value = 7537.25
$ws.Range("I18").Value = 1716.6666
$ws.Range("J18").Value = 24999
$ws.Range("K18").Value = 1716.6666
$ws.Range("L18").Value = 24999
$ws.Range("M18").Value = -1432.6666
$ws.Range("N18").Value = -25567

$ws.Range("H43").Value = 2725.6428
$ws.Range("I43").Value = 2842.0908
$ws.Range("J43").Value = 2298.6667
$ws.Range("K43").Value = 2842.0908
$ws.Range("L43").Value = 2298.6667
$ws.Range("M43").Value = -2773.0908
$ws.Range("N43").Value = -2436.6667

$ws.Range("H80").Value = 7408613.5
$ws.Range("I80").Value = 13889577
$ws.Range("J80").Value = 1798.7142
$ws.Range("K80").Value = 41668731
$ws.Range("L80").Value = 5396.142599999999
$ws.Range("M80").Value = -41667733
$ws.Range("N80").Value = -7392.142599999999

$ws.Range("H83").Value = 7408613.5
$ws.Range("I83").Value = 13889577
$ws.Range("J83").Value = 1798.7142
$ws.Range("K83").Value = 125006193
$ws.Range("L83").Value = 16188.4278
$ws.Range("M83").Value = -125001201
$ws.Range("N83").Value = -26172.4278

$ws.Range("H88").Value = 10717.315
$ws.Range("J88").Value = 12835.267
$ws.Range("L88").Value = 12835.267
$ws.Range("N88").Value = -13647.267

$ws.Range("H91").Value = 10717.315
$ws.Range("J91").Value = 12835.267
$ws.Range("L91").Value = 12835.267
$ws.Range("N91").Value = -15643.267

$ws.Range("H92").Value = 811.6667
$ws.Range("I92").Value = 296.85715
$ws.Range("K92").Value = 296.85715
$ws.Range("M92").Value = 951.14285

$ws.Range("H106").Value = 2349.611
$ws.Range("J106").Value = 2846.75
$ws.Range("L106").Value = 2846.75
$ws.Range("N106").Value = -4108.75

$ws.Range("H107").Value = 466.64285
$ws.Range("I107").Value = 422.33334
$ws.Range("K107").Value = 422.33334
$ws.Range("M107").Value = 1497.66666

$ws.Range("H138").Value = 1701.919
$ws.Range("I138").Value = 1050.9231
$ws.Range("K138").Value = 3152.7693
$ws.Range("M138").Value = 1987.2307

$ws.Range("H141").Value = 1556.25
$ws.Range("I141").Value = 1334.091
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 4002.273
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 1177.727
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4349.75
$ws.Range("J45").Value = 3133
$ws.Range("L45").Value = 3133
$ws.Range("N45").Value = -3887

$ws.Range("H61").Value = 6668076
$ws.Range("I61").Value = 6668076
$ws.Range("K61").Value = 6668076
$ws.Range("M61").Value = -6667864

$ws.Range("H88").Value = 1125.6875
$ws.Range("J88").Value = 1301
$ws.Range("L88").Value = 1301
$ws.Range("N88").Value = -2113

$ws.Range("H91").Value = 1125.6875
$ws.Range("J91").Value = 1301
$ws.Range("L91").Value = 1301
$ws.Range("N91").Value = -4109

$ws.Range("H136").Value = 6668076
$ws.Range("I136").Value = 6668076
$ws.Range("K136").Value = 20004228
$ws.Range("M136").Value = -20001678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6593.4287
$ws.Range("J80").Value = 3669
$ws.Range("L80").Value = 3669
$ws.Range("N80").Value = -5665

$ws.Range("H83").Value = 6593.4287
$ws.Range("J83").Value = 3669
$ws.Range("L83").Value = 18345
$ws.Range("N83").Value = -28329

$ws.Range("H86").Value = 5999.5
$ws.Range("I86").Value = 5999.5
$ws.Range("K86").Value = 5999.5
$ws.Range("M86").Value = -4876.5

$ws.Range("H89").Value = 5999.5
$ws.Range("I89").Value = 5999.5
$ws.Range("K89").Value = 29997.5
$ws.Range("M89").Value = -24381.5

$ws.Range("H94").Value = 1076.3529
$ws.Range("I94").Value = 837.61536
$ws.Range("J94").Value = 1852.25
$ws.Range("K94").Value = 837.61536
$ws.Range("L94").Value = 1852.25
$ws.Range("M94").Value = -386.61536
$ws.Range("N94").Value = -2754.25

$ws.Range("H134").Value = 2030727.9
$ws.Range("I134").Value = 2168403.2
$ws.Range("K134").Value = 6505209.600000001
$ws.Range("M134").Value = -6502674.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 837.2727
$ws.Range("I16").Value = 841.2857
$ws.Range("J16").Value = 830.25
$ws.Range("K16").Value = 841.2857
$ws.Range("L16").Value = 830.25
$ws.Range("M16").Value = -554.2857
$ws.Range("N16").Value = -1404.25

$ws.Range("H58").Value = 2060724.9
$ws.Range("I58").Value = 3087083.5
$ws.Range("K58").Value = 3087083.5
$ws.Range("M58").Value = -3086880.5

$ws.Range("H105").Value = 31347
$ws.Range("I105").Value = 36666.5
$ws.Range("J105").Value = 4749.5
$ws.Range("K105").Value = 36666.5
$ws.Range("L105").Value = 4749.5
$ws.Range("M105").Value = -34919.5
$ws.Range("N105").Value = -8243.5

$ws.Range("H113").Value = 837.2727
$ws.Range("I113").Value = 841.2857
$ws.Range("J113").Value = 830.25
$ws.Range("K113").Value = 841.2857
$ws.Range("L113").Value = 830.25
$ws.Range("M113").Value = 1328.7143
$ws.Range("N113").Value = -5170.25

$ws.Range("H122").Value = 4634.952
$ws.Range("I122").Value = 3752.3333
$ws.Range("J122").Value = 5811.778
$ws.Range("K122").Value = 11256.9999
$ws.Range("L122").Value = 17435.334
$ws.Range("M122").Value = -8806.999899999999
$ws.Range("N122").Value = -22335.334

$ws.Range("H134").Value = 7517.28
$ws.Range("I134").Value = 8223.727999999999
$ws.Range("J134").Value = 2336.6667
$ws.Range("K134").Value = 24671.184
$ws.Range("L134").Value = 7010.000100000001
$ws.Range("M134").Value = -22136.184
$ws.Range("N134").Value = -12080.0001

$ws.Range("H136").Value = 2060724.9
$ws.Range("I136").Value = 3087083.5
$ws.Range("K136").Value = 9261250.5
$ws.Range("M136").Value = -9258700.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 83.333336
$ws.Range("I41").Value = 50
$ws.Range("J41").Value = 150
$ws.Range("K41").Value = 150
$ws.Range("L41").Value = 450
$ws.Range("M41").Value = 188
$ws.Range("N41").Value = -1126

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4559.5557
$ws.Range("I40").Value = 4533.647
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4533.647
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4397.647
$ws.Range("N40").Value = -5272

$ws.Range("H132").Value = 1452920.2
$ws.Range("I132").Value = 2487293
$ws.Range("J132").Value = 4798.4
$ws.Range("K132").Value = 7461879
$ws.Range("L132").Value = 14395.2
$ws.Range("M132").Value = -7459349
$ws.Range("N132").Value = -19455.2

$ws.Range("H136").Value = 91423.42999999999
$ws.Range("J136").Value = 140547.56
$ws.Range("L136").Value = 421642.68
$ws.Range("N136").Value = -426742.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5257.2856
$ws.Range("I81").Value = 4950
$ws.Range("J81").Value = 5667
$ws.Range("K81").Value = 9900
$ws.Range("L81").Value = 11334
$ws.Range("M81").Value = -8839
$ws.Range("N81").Value = -13456

$ws.Range("H84").Value = 5257.2856
$ws.Range("I84").Value = 4950
$ws.Range("J84").Value = 5667
$ws.Range("K84").Value = 49500
$ws.Range("L84").Value = 56670
$ws.Range("M84").Value = -44196
$ws.Range("N84").Value = -67278

$ws.Range("H132").Value = 4794817
$ws.Range("I132").Value = 5298468.5
$ws.Range("K132").Value = 15895405.5
$ws.Range("M132").Value = -15892875.5

$ws.Range("H136").Value = 22494.8
$ws.Range("I136").Value = 28428.285
$ws.Range("K136").Value = 85284.855
$ws.Range("M136").Value = -82734.855
